$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 330.26315
$ws.Range("I5").Value = 126.5
$ws.Range("J5").Value = 3998
$ws.Range("K5").Value = 126.5
$ws.Range("L5").Value = 3998
$ws.Range("M5").Value = -11.5
$ws.Range("N5").Value = -4228
$ws.Range("H17").Value = 2166.8064
$ws.Range("J17").Value = 2166.8064
$ws.Range("L17").Value = 6500.4192
$ws.Range("N17").Value = -6836.4192
$ws.Range("H19").Value = 95
$ws.Range("I19").Value = 95
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 95
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 80
$ws.Range("N19").ClearContents()
$ws.Range("H40").Value = 5629.9
$ws.Range("I40").Value = 7001
$ws.Range("J40").Value = 5287.125
$ws.Range("K40").Value = 7001
$ws.Range("L40").Value = 5287.125
$ws.Range("M40").Value = -6826
$ws.Range("N40").Value = -5637.125
$ws.Range("H88").Value = 2468
$ws.Range("I88").Value = 1802
$ws.Range("J88").Value = 2601.2
$ws.Range("K88").Value = 1802
$ws.Range("L88").Value = 2601.2
$ws.Range("M88").Value = -1396
$ws.Range("N88").Value = -3413.2
$ws.Range("H91").Value = 2468
$ws.Range("I91").Value = 1802
$ws.Range("J91").Value = 2601.2
$ws.Range("K91").Value = 1802
$ws.Range("L91").Value = 2601.2
$ws.Range("M91").Value = -398
$ws.Range("N91").Value = -5409.2
$ws.Range("H97").Value = 2566.182
$ws.Range("I97").Value = 999
$ws.Range("J97").Value = 2722.9
$ws.Range("K97").Value = 2997
$ws.Range("L97").Value = 8168.700000000001
$ws.Range("N97").Value = -9160.700000000001
$ws.Range("M97").Value = -2501
$ws.Range("H100").Value = 2971.2856
$ws.Range("I100").Value = 2825
$ws.Range("K100").Value = 2825
$ws.Range("M100").Value = -2284
$ws.Range("H107").Value = 40573.4
$ws.Range("I107").Value = 40573.4
$ws.Range("K107").Value = 40573.4
$ws.Range("M107").Value = -38653.4
$ws.Range("H123").Value = 70766.8
$ws.Range("J123").Value = 70766.8
$ws.Range("L123").Value = 70766.8
$ws.Range("N123").Value = -80566.8
$ws.Range("H125").Value = 10105413
$ws.Range("J125").Value = 12350571
$ws.Range("L125").Value = 111155139
$ws.Range("N125").Value = -111160059
$ws.Range("H131").Value = 5325
$ws.Range("I131").Value = 2870.6667
$ws.Range("J131").Value = 7333.091
$ws.Range("K131").Value = 8612.000100000001
$ws.Range("L131").Value = 21999.273
$ws.Range("M131").Value = -3572.000100000001
$ws.Range("N131").Value = -32079.273
$ws.Range("H135").Value = 770669.9399999999
$ws.Range("I135").Value = 1001217.56
$ws.Range("K135").Value = 9010958.040000001
$ws.Range("M135").Value = -9008423.040000001
$ws.Range("H137").Value = 2870.6086
$ws.Range("I137").Value = 2695
$ws.Range("J137").Value = 3272
$ws.Range("K137").Value = 8085
$ws.Range("L137").Value = 9816
$ws.Range("M137").Value = -5535
$ws.Range("N137").Value = -14916

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2357.3684
$ws.Range("I61").Value = 2181.8235
$ws.Range("K61").Value = 2181.8235
$ws.Range("M61").Value = -1969.8235
$ws.Range("H63").Value = 5555.6665
$ws.Range("J63").Value = 7389
$ws.Range("L63").Value = 7389
$ws.Range("N63").Value = -8761
$ws.Range("H66").Value = 5555.6665
$ws.Range("J66").Value = 7389
$ws.Range("L66").Value = 36945
$ws.Range("N66").Value = -43809
$ws.Range("H74").Value = 1612.3334
$ws.Range("I74").Value = 1612.3334
$ws.Range("K74").Value = 1612.3334
$ws.Range("M74").Value = -738.3334
$ws.Range("H77").Value = 1612.3334
$ws.Range("I77").Value = 1612.3334
$ws.Range("K77").Value = 8061.666999999999
$ws.Range("M77").Value = -3693.666999999999
$ws.Range("H102").Value = 1793.1666
$ws.Range("I102").Value = 1762
$ws.Range("K102").Value = 1762
$ws.Range("M102").Value = -140
$ws.Range("H110").Value = 359555
$ws.Range("J110").Value = 1100
$ws.Range("L110").Value = 1100
$ws.Range("N110").Value = -5190
$ws.Range("H132").Value = 3034.5757
$ws.Range("I132").Value = 3347.12
$ws.Range("J132").Value = 2057.875
$ws.Range("K132").Value = 10041.36
$ws.Range("L132").Value = 6173.625
$ws.Range("M132").Value = -7511.360000000001
$ws.Range("N132").Value = -11233.625
$ws.Range("H136").Value = 2357.3684
$ws.Range("I136").Value = 2181.8235
$ws.Range("K136").Value = 6545.470499999999
$ws.Range("M136").Value = -3995.470499999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 107499.5
$ws.Range("J35").Value = 107499.5
$ws.Range("L35").Value = 107499.5
$ws.Range("N35").Value = -108119.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 2265.6155
$ws.Range("J5").Value = 2529.5715
$ws.Range("L5").Value = 2529.5715
$ws.Range("N5").Value = -2753.5715
$ws.Range("H7").Value = 517.2083
$ws.Range("I7").Value = 542.2778
$ws.Range("J7").Value = 442
$ws.Range("K7").Value = 542.2778
$ws.Range("L7").Value = 442
$ws.Range("M7").Value = -429.2778
$ws.Range("N7").Value = -668
$ws.Range("H107").Value = 554.8333
$ws.Range("J107").Value = 909.8
$ws.Range("L107").Value = 909.8
$ws.Range("N107").Value = -4749.8
$ws.Range("H132").Value = 1435
$ws.Range("I132").Value = 1196.8
$ws.Range("J132").Value = 1832
$ws.Range("K132").Value = 3590.4
$ws.Range("L132").Value = 5496
$ws.Range("M132").Value = -1060.4
$ws.Range("N132").Value = -10556
$ws.Range("H134").Value = 214820.6
$ws.Range("I134").Value = 2130.4
$ws.Range("J134").Value = 5000350
$ws.Range("K134").Value = 6391.200000000001
$ws.Range("L134").Value = 15001050
$ws.Range("M134").Value = -3856.200000000001
$ws.Range("N134").Value = -15006120
$ws.Range("H141").Value = 435869.56
$ws.Range("J141").Value = 488465.16
$ws.Range("L141").Value = 488465.16
$ws.Range("N141").Value = -498825.16

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 385694.94
$ws.Range("I5").Value = 980.63635
$ws.Range("J5").Value = 667818.75
$ws.Range("K5").Value = 2941.90905
$ws.Range("L5").Value = 2003456.25
$ws.Range("M5").Value = -2829.90905
$ws.Range("N5").Value = -2003680.25
$ws.Range("H33").Value = 4226444
$ws.Range("I33").Value = 12345853
$ws.Range("J33").Value = 166739.17
$ws.Range("K33").Value = 74075118
$ws.Range("L33").Value = 1000435.02
$ws.Range("M33").Value = -74074835
$ws.Range("N33").Value = -1001001.02
$ws.Range("H113").Value = 3087444.5
$ws.Range("J113").Value = 985.7143
$ws.Range("L113").Value = 2957.1429
$ws.Range("N113").Value = -7297.1429
$ws.Range("H118").Value = 961.4286
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H132").Value = 737241.5600000001
$ws.Range("I132").Value = 143802.86
$ws.Range("J132").Value = 1256500.4
$ws.Range("K132").Value = 1294225.74
$ws.Range("L132").Value = 11308503.6
$ws.Range("M132").Value = -1291695.74
$ws.Range("N132").Value = -11313563.6
$ws.Range("H135").Value = 385694.94
$ws.Range("I135").Value = 980.63635
$ws.Range("J135").Value = 667818.75
$ws.Range("K135").Value = 8825.727150000001
$ws.Range("L135").Value = 6010368.75
$ws.Range("M135").Value = -6290.727150000001
$ws.Range("N135").Value = -6015438.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 25000
$ws.Range("I97").Value = 25000
$ws.Range("K97").Value = 25000
$ws.Range("M97").Value = -24504
$ws.Range("H122").Value = 4509.8887
$ws.Range("I122").Value = 2563
$ws.Range("K122").Value = 7689
$ws.Range("M122").Value = -5239

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 129499
$ws.Range("I40").Value = 202398.4
$ws.Range("K40").Value = 202398.4
$ws.Range("M40").Value = -202262.4
$ws.Range("H46").Value = 4307.615
$ws.Range("J46").Value = 5583.3335
$ws.Range("L46").Value = 5583.3335
$ws.Range("N46").Value = -5959.3335
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H132").Value = 2966.6667
$ws.Range("J132").Value = 3001
$ws.Range("L132").Value = 9003
$ws.Range("N132").Value = -14063

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 52635356
$ws.Range("J122").Value = 6001.25
$ws.Range("L122").Value = 18003.75
$ws.Range("N122").Value = -22903.75
$ws.Range("H124").Value = 84981.336
$ws.Range("J124").Value = 84981.336
$ws.Range("L124").Value = 84981.336
$ws.Range("N124").Value = -94801.336
$ws.Range("H132").Value = 147415
$ws.Range("H133").Value = 63999.5
$ws.Range("J133").Value = 63999.5
$ws.Range("L133").Value = 63999.5
$ws.Range("N133").Value = -74119.5
$ws.Range("H136").Value = 9089854
$ws.Range("I136").Value = 11446765
$ws.Range("K136").Value = 34340295
$ws.Range("M136").Value = -34337745
